# "Ran Prod Verification Script"
# The Katalon "UM-Data-Prod" bootstrap workbook stamps the DateProd cell
# (column B, row 2) of every test-data sheet with the timestamp of the last
# Prod verification run. Re-running the script against Prod updates that
# single cell on each sheet.

$wb = $excel.ActiveWorkbook

$runTimes = @{
    "CreateUser"    = "Thu Aug 28 07:31:00 IST 2025"
    "FindUser"      = "Thu Aug 28 07:31:41 IST 2025"
    "ModifyUser"    = "Thu Aug 28 07:32:18 IST 2025"
    "ModifyUserPwd" = "Thu Aug 28 07:33:04 IST 2025"
    "FindCaseUser"  = "Thu Aug 28 07:34:05 IST 2025"
    "AddDeleteRole" = "Thu Aug 28 07:29:53 IST 2025"
    "SearchRole"    = "Thu Aug 28 07:30:31 IST 2025"
}

foreach ($sheetName in $runTimes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B2").Value = $runTimes[$sheetName]
}
